# Bai 24.4 (PageRank) - fix the random-jump term so the inequality reads
# "... khong nho hon alpha/N voi alpha la xac suat nhay ngau nhien."
# i.e. insert "/N" right after the first standalone "alpha" that is
# followed by " voi alpha la xac suat...".

$p = $ppt.ActivePresentation

$needle = "hơn α với"
$targetSlide = $null
$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text.IndexOf($needle) -ge 0) {
                $targetSlide = $slide
                $targetShape = $shape
            }
        }
    }
}

if ($targetShape -eq $null) {
    throw "Could not find shape containing '$needle'"
}

$tr = $targetShape.TextFrame.TextRange
$full = $tr.Text
$idx = $full.IndexOf($needle)

# $needle = "hơn α với" -> 'h' 'ơ' 'n' ' ' 'α' ' ' 'v' 'ơ' 'i'
#             index:        0   1   2   3   4   5 ...
# The space right after "α" sits at offset 5 within the needle.
$spacePos0 = $idx + 5
$spaceRange = $tr.Characters($spacePos0 + 1, 1)

if ($spaceRange.Text -ne " ") {
    throw "Unexpected character at target position: [$($spaceRange.Text)]"
}

$spaceRange.Text = "/N "
